$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.215.30'
$ws.Range('E2').Value = '  -0.97%  '

$ws.Range('D3').Value = '2.304.24'
$ws.Range('E3').Value = '  -2.97%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = "'309.21"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -6.21%  '

$ws.Range('D6').Value = "'105.30"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.94%  '

$ws.Range('D7').Value = "'0.624"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.00%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = "'0.603"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.74%  '

$ws.Range('D10').Value = "'39.62"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.44%  '

$ws.Range('D11').Value = "'0.0909"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.99%  '

$ws.Range('D12').Value = "'8.25"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.60%  '

$ws.Range('D13').Value = "'0.105"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.03%  '

$ws.Range('D14').Value = "'0.964"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.33%  '

$ws.Range('D15').Value = "'15.34"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.94%  '

$ws.Range('D16').Value = '2.660.61'
$ws.Range('E16').Value = '  -2.91%  '

$ws.Range('D17').Value = '2.324.55'
$ws.Range('E17').Value = '  -1.75%  '

$ws.Range('D18').Value = '42.202.24'
$ws.Range('E18').Value = '  -1.03%  '

$ws.Range('D19').Value = "'7.40"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.78%  '

$ws.Range('D20').Value = "'0.0000104"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.06%  '

$ws.Range('D21').Value = "'74.82"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.34%  '

$ws.Range('D22').Value = "'3.44"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -8.43%  '

$ws.Range('D23').Value = "'258.95"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.16%  '

$ws.Range('D24').Value = "'2.26"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.02%  '

$ws.Range('D25').Value = "'9.15"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.67%  '

$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.41%  '

$ws.Range('D27').Value = "'10.95"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.14%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = "'2.26"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.68%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = "'22.90"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.29%  '

$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = "'165.04"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.72%  '

$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = "'35.58"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.04%  '

$ws.Range('D32').Value = "'0.0887"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.14%  '

$ws.Range('E33').Value = '  -6.51%  '

$ws.Range('D34').Value = "'5.83"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.00%  '

$ws.Range('D35').Value = "'0.117"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +10.71%  '

$ws.Range('D36').Value = "'0.128"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.28%  '

$ws.Range('D37').Value = "'4.45"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.82%  '

$ws.Range('D38').Value = "'0.0348"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.60%  '

$ws.Range('D39').Value = "'3.66"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.40%  '

$ws.Range('D40').Value = "'2.62"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.93%  '

$ws.Range('D41').Value = "'100.05"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.00%  '

$ws.Range('D42').Value = "'1.46"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.58%  '

$ws.Range('D43').Value = "'69.63"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.26%  '

$ws.Range('D44').Value = "'0.229"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.12%  '

$ws.Range('E45').Value = '  +0.16%  '

$ws.Range('D46').Value = "'11.99"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.34%  '

$ws.Range('D47').Value = "'110.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.69%  '

$ws.Range('D48').Value = "'5.38"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.84%  '

$ws.Range('D49').Value = "'8.97"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.94%  '

$ws.Range('D50').Value = "'73.38"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.00%  '

$ws.Range('D51').Value = "'1.25"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.26%  '
